# Update TPM values for Wnt10b-Fzd7 LR pairs, and add Neutrophils sending-cluster rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt10b"
$ws.Range("C2").Value = "Fzd7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.044796
$ws.Range("H2").Value = 0.134388
$ws.Range("I2").Value = 0.2628163758609765
$ws.Range("J2").Value = 0.3020769645747168
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.905108
$ws.Range("N2").Value = 3.810216
$ws.Range("O2").Value = 0.07580486173280727
$ws.Range("P2").Value = 0.05715529216076502
$ws.Range("Q2").Value = 0.08534121796800001
$ws.Range("R2").Value = 0.512047307808
$ws.Range("S2").Value = 0.01992275903325883
$ws.Range("T2").Value = 0.017265297165305

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt10b"
$ws.Range("C3").Value = "Fzd7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.044796
$ws.Range("H3").Value = 0.134388
$ws.Range("I3").Value = 0.2628163758609765
$ws.Range("J3").Value = 0.3020769645747168
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.068283333333333
$ws.Range("N3").Value = 9.20485
$ws.Range("O3").Value = 0.1220879833796353
$ws.Range("P3").Value = 0.1380777076800943
$ws.Range("Q3").Value = 0.1374468202
$ws.Range("R3").Value = 1.2370213818
$ws.Range("S3").Value = 0.03208672132801089
$ws.Range("T3").Value = 0.04171009481143793

# Row 4: FAPs -> Inflammatory-Mac
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt10b"
$ws.Range("C4").Value = "Fzd7"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.044796
$ws.Range("H4").Value = 0.134388
$ws.Range("I4").Value = 0.2628163758609765
$ws.Range("J4").Value = 0.3020769645747168
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.725954
$ws.Range("N4").Value = 11.177862
$ws.Range("O4").Value = 0.148256911310435
$ws.Range("P4").Value = 0.1676739503331867
$ws.Range("Q4").Value = 0.166907835384
$ws.Range("R4").Value = 1.502170518456
$ws.Range("S4").Value = 0.03896434412695074
$ws.Range("T4").Value = 0.05065043795490085

# Row 5: FAPs -> MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt10b"
$ws.Range("C5").Value = "Fzd7"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.044796
$ws.Range("H5").Value = 0.134388
$ws.Range("I5").Value = 0.2628163758609765
$ws.Range("J5").Value = 0.3020769645747168
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.825836
$ws.Range("N5").Value = 13.651672
$ws.Range("O5").Value = 0.2716022158275637
$ws.Range("P5").Value = 0.2047824327132465
$ws.Range("Q5").Value = 0.305770149456
$ws.Range("R5").Value = 1.834620896736
$ws.Range("S5").Value = 0.07138151003961105
$ws.Range("T5").Value = 0.0618600556722437

# Row 6: FAPs -> Neutrophils
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt10b"
$ws.Range("C6").Value = "Fzd7"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.044796
$ws.Range("H6").Value = 0.134388
$ws.Range("I6").Value = 0.2628163758609765
$ws.Range("J6").Value = 0.3020769645747168
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.247626666666666
$ws.Range("N6").Value = 12.74288
$ws.Range("O6").Value = 0.1690144349607748
$ws.Range("P6").Value = 0.1911500632430207
$ws.Range("Q6").Value = 0.19027668416
$ws.Range("R6").Value = 1.71249015744
$ws.Range("S6").Value = 0.04441976126458155
$ws.Range("T6").Value = 0.05774203088271684

# Row 7: FAPs -> Resolving-Mac
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt10b"
$ws.Range("C7").Value = "Fzd7"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.044796
$ws.Range("H7").Value = 0.134388
$ws.Range("I7").Value = 0.2628163758609765
$ws.Range("J7").Value = 0.3020769645747168
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.358931000000001
$ws.Range("N7").Value = 16.076793
$ws.Range("O7").Value = 0.213233592788784
$ws.Range("P7").Value = 0.2411605538696867
$ws.Range("Q7").Value = 0.2400586730760001
$ws.Range("R7").Value = 2.160528057684
$ws.Range("S7").Value = 0.05604128006856346
$ws.Range("T7").Value = 0.07284904808811242

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Wnt10b"
$ws.Range("C8").Value = "Fzd7"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.066458
$ws.Range("H8").Value = 0.132916
$ws.Range("I8").Value = 0.3899064806448963
$ws.Range("J8").Value = 0.2987682071569862
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.905108
$ws.Range("N8").Value = 3.810216
$ws.Range("O8").Value = 0.07580486173280727
$ws.Range("P8").Value = 0.05715529216076502
$ws.Range("Q8").Value = 0.126609667464
$ws.Range("R8").Value = 0.5064386698560001
$ws.Range("S8").Value = 0.02955680685401186
$ws.Range("T8").Value = 0.01707618416840551

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Wnt10b"
$ws.Range("C9").Value = "Fzd7"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.066458
$ws.Range("H9").Value = 0.132916
$ws.Range("I9").Value = 0.3899064806448963
$ws.Range("J9").Value = 0.2987682071569862
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.068283333333333
$ws.Range("N9").Value = 9.20485
$ws.Range("O9").Value = 0.1220879833796353
$ws.Range("P9").Value = 0.1380777076800943
$ws.Range("Q9").Value = 0.2039119737666667
$ws.Range("R9").Value = 1.2234718426
$ws.Range("S9").Value = 0.04760289592858621
$ws.Range("T9").Value = 0.04125322917192818

# Row 10: MuSCs -> Inflammatory-Mac
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Wnt10b"
$ws.Range("C10").Value = "Fzd7"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.066458
$ws.Range("H10").Value = 0.132916
$ws.Range("I10").Value = 0.3899064806448963
$ws.Range("J10").Value = 0.2987682071569862
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.725954
$ws.Range("N10").Value = 11.177862
$ws.Range("O10").Value = 0.148256911310435
$ws.Range("P10").Value = 0.1676739503331867
$ws.Range("Q10").Value = 0.247619450932
$ws.Range("R10").Value = 1.485716705592
$ws.Range("S10").Value = 0.05780633052033422
$ws.Range("T10").Value = 0.05009564552797573

# Row 11: MuSCs -> MuSCs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Wnt10b"
$ws.Range("C11").Value = "Fzd7"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.066458
$ws.Range("H11").Value = 0.132916
$ws.Range("I11").Value = 0.3899064806448963
$ws.Range("J11").Value = 0.2987682071569862
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 6.825836
$ws.Range("N11").Value = 13.651672
$ws.Range("O11").Value = 0.2716022158275637
$ws.Range("P11").Value = 0.2047824327132465
$ws.Range("Q11").Value = 0.453631408888
$ws.Range("R11").Value = 1.814525635552
$ws.Range("S11").Value = 0.1058994641086809
$ws.Range("T11").Value = 0.06118248027898282

# Row 12: MuSCs -> Neutrophils
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Wnt10b"
$ws.Range("C12").Value = "Fzd7"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.066458
$ws.Range("H12").Value = 0.132916
$ws.Range("I12").Value = 0.3899064806448963
$ws.Range("J12").Value = 0.2987682071569862
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.247626666666666
$ws.Range("N12").Value = 12.74288
$ws.Range("O12").Value = 0.1690144349607748
$ws.Range("P12").Value = 0.1911500632430207
$ws.Range("Q12").Value = 0.2822887730133333
$ws.Range("R12").Value = 1.69373263808
$ws.Range("S12").Value = 0.0658998235137414
$ws.Range("T12").Value = 0.05710956169306181

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Wnt10b"
$ws.Range("C13").Value = "Fzd7"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.066458
$ws.Range("H13").Value = 0.132916
$ws.Range("I13").Value = 0.3899064806448963
$ws.Range("J13").Value = 0.2987682071569862
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 5.358931000000001
$ws.Range("N13").Value = 16.076793
$ws.Range("O13").Value = 0.213233592788784
$ws.Range("P13").Value = 0.2411605538696867
$ws.Range("Q13").Value = 0.3561438363980001
$ws.Range("R13").Value = 2.136863018388
$ws.Range("S13").Value = 0.0831411597195417
$ws.Range("T13").Value = 0.07205110631663207

# Row 14: Neutrophils -> ECs
$ws.Range("A14").Value = "Neutrophils"
$ws.Range("B14").Value = "Wnt10b"
$ws.Range("C14").Value = "Fzd7"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.059192
$ws.Range("H14").Value = 0.177576
$ws.Range("I14").Value = 0.3472771434941271
$ws.Range("J14").Value = 0.3991548282682971
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.905108
$ws.Range("N14").Value = 3.810216
$ws.Range("O14").Value = 0.07580486173280727
$ws.Range("P14").Value = 0.05715529216076502
$ws.Range("Q14").Value = 0.112767152736
$ws.Range("R14").Value = 0.676602916416
$ws.Range("S14").Value = 0.02632529584553658
$ws.Range("T14").Value = 0.02281381082705451

# Row 15: Neutrophils -> FAPs
$ws.Range("A15").Value = "Neutrophils"
$ws.Range("B15").Value = "Wnt10b"
$ws.Range("C15").Value = "Fzd7"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.059192
$ws.Range("H15").Value = 0.177576
$ws.Range("I15").Value = 0.3472771434941271
$ws.Range("J15").Value = 0.3991548282682971
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.068283333333333
$ws.Range("N15").Value = 9.20485
$ws.Range("O15").Value = 0.1220879833796353
$ws.Range("P15").Value = 0.1380777076800943
$ws.Range("Q15").Value = 0.1816178270666667
$ws.Range("R15").Value = 1.6345604436
$ws.Range("S15").Value = 0.04239836612303823
$ws.Range("T15").Value = 0.05511438369672815

# Row 16: Neutrophils -> Inflammatory-Mac
$ws.Range("A16").Value = "Neutrophils"
$ws.Range("B16").Value = "Wnt10b"
$ws.Range("C16").Value = "Fzd7"
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.059192
$ws.Range("H16").Value = 0.177576
$ws.Range("I16").Value = 0.3472771434941271
$ws.Range("J16").Value = 0.3991548282682971
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 3.725954
$ws.Range("N16").Value = 11.177862
$ws.Range("O16").Value = 0.148256911310435
$ws.Range("P16").Value = 0.1676739503331867
$ws.Range("Q16").Value = 0.220546669168
$ws.Range("R16").Value = 1.984920022512
$ws.Range("S16").Value = 0.05148623666315
$ws.Range("T16").Value = 0.06692786685031012

# Row 17: Neutrophils -> MuSCs
$ws.Range("A17").Value = "Neutrophils"
$ws.Range("B17").Value = "Wnt10b"
$ws.Range("C17").Value = "Fzd7"
$ws.Range("D17").Value = "MuSCs"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.059192
$ws.Range("H17").Value = 0.177576
$ws.Range("I17").Value = 0.3472771434941271
$ws.Range("J17").Value = 0.3991548282682971
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 6.825836
$ws.Range("N17").Value = 13.651672
$ws.Range("O17").Value = 0.2716022158275637
$ws.Range("P17").Value = 0.2047824327132465
$ws.Range("Q17").Value = 0.404034884512
$ws.Range("R17").Value = 2.424209307072
$ws.Range("S17").Value = 0.09432124167927174
$ws.Range("T17").Value = 0.08173989676202004

# Row 18: Neutrophils -> Neutrophils
$ws.Range("A18").Value = "Neutrophils"
$ws.Range("B18").Value = "Wnt10b"
$ws.Range("C18").Value = "Fzd7"
$ws.Range("D18").Value = "Neutrophils"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.059192
$ws.Range("H18").Value = 0.177576
$ws.Range("I18").Value = 0.3472771434941271
$ws.Range("J18").Value = 0.3991548282682971
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 4.247626666666666
$ws.Range("N18").Value = 12.74288
$ws.Range("O18").Value = 0.1690144349607748
$ws.Range("P18").Value = 0.1911500632430207
$ws.Range("Q18").Value = 0.2514255176533333
$ws.Range("R18").Value = 2.26282965888
$ws.Range("S18").Value = 0.0586948501824518
$ws.Range("T18").Value = 0.07629847066724206

# Row 19: Neutrophils -> Resolving-Mac
$ws.Range("A19").Value = "Neutrophils"
$ws.Range("B19").Value = "Wnt10b"
$ws.Range("C19").Value = "Fzd7"
$ws.Range("D19").Value = "Resolving-Mac"
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.059192
$ws.Range("H19").Value = 0.177576
$ws.Range("I19").Value = 0.3472771434941271
$ws.Range("J19").Value = 0.3991548282682971
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 5.358931000000001
$ws.Range("N19").Value = 16.076793
$ws.Range("O19").Value = 0.213233592788784
$ws.Range("P19").Value = 0.2411605538696867
$ws.Range("Q19").Value = 0.3172058437520001
$ws.Range("R19").Value = 2.854852593768001
$ws.Range("S19").Value = 0.07405115300067881
$ws.Range("T19").Value = 0.0962603994649422
